$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells (row 1)
$ws.Range("C1").Value = "Stock_num"
$ws.Range("D1").Value = "PrecioEnDolares_num"
$ws.Range("E1").Value = "PrecioEnpesos_num"
$ws.Range("G1").Value = "Fecha_Ingreso_date"

# 2. Row 3: convert D3 from text "699.00" to the number 699
$ws.Range("D3").Value = 699

# 3. Row 4: a new, empty styled cell appears at J4 (matches the style already
#    used by the stray empty cell at I30)
$ws.Range("I30").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Row 5: D5 becomes the text "hola"
$ws.Range("D5").Value = "hola"

# 5. Row 6: G6 becomes the text "hola" (was a date)
$ws.Range("G6").Value = "hola"

# 6. Row 24: D24 becomes the text "hola"
$ws.Range("D24").Value = "hola"

# 7. Row 25: D25 keeps its text ("149.00") but takes on the special
#    underline/white-on-dark style already used at F19/D38
$ws.Range("F19").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 8. Column width tweaks
$ws.Columns("D").ColumnWidth = 21.5703125
$ws.Columns("E").ColumnWidth = 25.140625
$ws.Columns("G").ColumnWidth = 22.140625

# 9. Update the view: scroll position and active selection
$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 20
